$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = "ДГИ-В-87616/24  14.08.2024"
$ws.Range("G2").Value = "ДГИ-В-87616/24 14.08.2024"
$ws.Range("H2").Value = "Кому: Мусиенко О.А. (Департамент городского имущества города Москвы)  От кого:  Демонова Л.В. (Департамент городского имущества города Москвы)"
$ws.Range("I2").Value = "О показе жилых помещений ЮВАО"

# Row 3
$ws.Range("E3").Value = "Вн"
$ws.Range("F3").Value = "ДГИ-В-87613/24  14.08.2024"
$ws.Range("G3").Value = "ДГИ-В-87613/24 14.08.2024"
$ws.Range("H3").Value = "Кому: Мусиенко О.А. (Департамент городского имущества города Москвы)  От кого:  Демонова Л.В. (Департамент городского имущества города Москвы)"
$ws.Range("I3").Value = "Об осмотре жилых помещений по КПИ и ДСН в ЮВАО"

# Row 4
$ws.Range("F4").Value = "ДГИ-1-44675/24  13.08.2024"
$ws.Range("G4").Value = "4 02.08.2024"
$ws.Range("H4").Value = "Кому: Гаман М.Ф. (Департамент городского имущества города Москвы)  От кого:  Ермоленко Н.В. (Финансовый управляющий)"
$ws.Range("I4").Value = "ДГИ-241024/24-(0)-0 запрос по делу А40-95663/2024"

# Row 5
$ws.Range("F5").Value = "ДГИ-1-44670/24  13.08.2024"
$ws.Range("G5").Value = "2 05.08.2024"
$ws.Range("H5").Value = "Кому: Гаман М.Ф. (Департамент городского имущества города Москвы)  От кого:  Мокрушин С.В. (Конкурсный управляющий)"
$ws.Range("I5").Value = "ДГИ-241090/24-(0)-0 запрос по делу А40-287374/23"

# Row 6
$ws.Range("F6").Value = "ДГИ-1-44664/24  13.08.2024"
$ws.Range("G6").Value = "Б/Н 02.08.2024"
$ws.Range("H6").Value = "Кому: Гаман М.Ф. (Департамент городского имущества города Москвы)  От кого:  Алехин Н.Н. (Конкурсный управляющий)"
$ws.Range("I6").Value = "ДГИ-240927/24-(0)-0 запрос по делу А40-169117/23"

# Row 7
$ws.Range("F7").Value = "ДГИ-1-44660/24  13.08.2024"
$ws.Range("G7").Value = "9 02.08.2024"
$ws.Range("H7").Value = "Кому: Гаман М.Ф. (Департамент городского имущества города Москвы)  От кого:  Курзин Д.А. (Финансовый управляющий)"
$ws.Range("I7").Value = "ДГИ-240954/24-(0)-0 запрос по делу А40-111870/2024"

# Row 8
$ws.Range("F8").Value = "ДГИ-1-44658/24  13.08.2024"
$ws.Range("G8").Value = "1 04.06.2024"
$ws.Range("H8").Value = "Кому: Гаман М.Ф. (Департамент городского имущества города Москвы)  От кого:  Вахрушев В.О. (Временный управляющий)"
$ws.Range("I8").Value = "ДГИ-240946/24-(0)-0 запрос по делу А40-79798/24"

# Row 9
$ws.Range("F9").Value = "ДГИ-1-44654/24  13.08.2024"
$ws.Range("G9").Value = "9 05.08.2024"
$ws.Range("H9").Value = "Кому: Гаман М.Ф. (Департамент городского имущества города Москвы)  От кого:  Ломакина М.М. (Финансовый управляющий)"
$ws.Range("I9").Value = "ДГИ-241014/24-(0)-0 запрос по делу А40-73286/2024"

# Row 10
$ws.Range("F10").Value = "ДГИ-1-44650/24  13.08.2024"
$ws.Range("G10").Value = "859-47 07.08.2024"
$ws.Range("H10").Value = "Кому: Гаман М.Ф. (Департамент городского имущества города Москвы)  От кого:  Стародубцев А.В. (Конкурсный управляющий)"
$ws.Range("I10").Value = "ДГИ-240940/24-(0)-0 запрос по делу А40-253586/23"

# Row 11
$ws.Range("E11").Value = "Гр"
$ws.Range("F11").Value = "ДГИ-ЭГР-46848/24  13.08.2024"
$ws.Range("G11").Value = "56021949 13.08.2024"
$ws.Range("H11").Value = "Кому: Гаман М.Ф. (Департамент городского имущества города Москвы)  От кого:  Обращение граждан (Обращение граждан)"
$ws.Range("I11").Value = "Обращения граждан Вопрос 1. Сообщение с mos.ru, идентификатор: 56021949 Корякин Анатолий Алексеевич, по очереди"

# Row 12
$ws.Range("F12").Value = "ДГИ-1-44633/24  13.08.2024"
$ws.Range("G12").Value = "20 06.08.2024"
$ws.Range("H12").Value = "Кому: Гаман М.Ф. (Департамент городского имущества города Москвы)  От кого:  Пуртов Н.С. (Финансовый управляющий)"
$ws.Range("I12").Value = "ДГИ-240886/24-(0)-0 запрос по делу А40-284028/23"

# Row 13
$ws.Range("F13").Value = "ДГИ-Э-128678/24  13.08.2024"
$ws.Range("G13").Value = "А40-137492/2024 13.08.2024"
$ws.Range("H13").Value = "Кому: Гаман М.Ф. (Департамент городского имущества города Москвы)  От кого:  Волкова Л.И. (Финансовый управляющий)"
$ws.Range("I13").Value = "Запрос по делу №А40-137492/2024"

# Row 14
$ws.Range("F14").Value = "ДГИ-1-44632/24  13.08.2024"
$ws.Range("G14").Value = "27 05.08.2024"
$ws.Range("H14").Value = "Кому: Гаман М.Ф. (Департамент городского имущества города Москвы)  От кого:  Феденко А.А. (Финансовый управляющий)"
$ws.Range("I14").Value = "ДГИ-240884/24-(0)-0 запрос по делу А40-103528/23"

# Row 15
$ws.Range("F15").Value = "ДГИ-1-44630/24  13.08.2024"
$ws.Range("G15").Value = "3 29.07.2024"
$ws.Range("H15").Value = "Кому: Гаман М.Ф. (Департамент городского имущества города Москвы)  От кого:  Слесарев К.И. (Финансовый управляющий)"
$ws.Range("I15").Value = "ДГИ-240883/24-(0)-0 запрос по делу А40-116370/24"

# Row 16
$ws.Range("E16").Value = "Гр"
$ws.Range("F16").Value = "ДГИ-ЭГР-46845/24  13.08.2024"
$ws.Range("G16").Value = "56021855 13.08.2024"
$ws.Range("H16").Value = "Кому: Гаман М.Ф. (Департамент городского имущества города Москвы)  От кого:  Обращение граждан (Обращение граждан)"
$ws.Range("I16").Value = "Обращения граждан Вопрос 1. Сообщение с mos.ru, идентификатор: 56021855 Фоломеева Наталья Ивановна, Отказ в принятии на учет нуждающихся в жилых помещениях"
